# "add loger and logic for search in alternative region"
#
# The tracked sheet ("выгрузка_PU") held a list of case numbers in column A
# (with an occasional "нет" flag in column B). The commit trims the sheet
# down to a single search number and drops the now-unused flag column:
#   - A2 becomes the new search number "07/048574" (was "02/056719")
#   - rows 3:15 (the rest of the old number list) are removed
#   - column B (the "нет" flags) is removed entirely
#   - the active selection moves to A2

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New alternative-region search number.
$ws.Range("A2").Value = "07/048574"

# Drop the remaining old rows of numbers/flags (rows 3-15); row 16 (already
# blank) shifts up to become the new, empty row 3.
$ws.Rows("3:15").Delete()

# The "нет" flag column is no longer used - remove it.
$ws.Columns("B").Delete()

# Match the committed selection state (active cell A2).
$ws.Range("A2").Select() | Out-Null
